# Generate Report for handback
# Adds a new handback record (row 4) for file
# d0e42c9f-4808-4a5f-ae20-7dae63366bda to the Overview, zh-cn and de-de
# worksheets of the handback-status report, mirroring the existing rows.

$wb = $excel.ActiveWorkbook

$guid   = "d0e42c9f-4808-4a5f-ae20-7dae63366bda"
$xlfSum = "59c6eead89ddadf9b36ddd18738065621f482cc4"

$mdName     = "$guid.md"
$zhXlfName  = "$guid.$xlfSum.zh-cn.xlf"
$deXlfName  = "$guid.$xlfSum.de-de.xlf"

$statusInSync = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: new row 4 -> File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/60aed80794e1226e8469026e0ce8cf0cadddc531/e2e/$mdName",
    $null,
    $null,
    $mdName)
$wsOverview.Range("B4").Value = $statusInSync
$wsOverview.Range("C4").Value = $statusInSync

# ---------------------------------------------------------------------
# zh-cn sheet: new row 4
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/60aed80794e1226e8469026e0ce8cf0cadddc531/e2e/$mdName",
    $null,
    $null,
    $mdName)
$wsZhCn.Range("B4").Value = $statusInSync
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e29a7f0430dd0f84f4fa387d69f2623558b3a4bd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlfName",
    $null,
    $null,
    $zhXlfName)
$wsZhCn.Range("D4").Value = "2016-01-25 10:47:39"
$wsZhCn.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/2766952e235a1b3db63e6e592972ca2d5705d76e/e2e/$mdName",
    $null,
    $null,
    $mdName)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/decdd909402c5f3c4ee2404aea33fdcf098e7c00/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlfName",
    $null,
    $null,
    $zhXlfName)
$wsZhCn.Range("G4").Value = "2016-01-25 10:48:19"
$wsZhCn.Range("H4").Value = "Include"

# ---------------------------------------------------------------------
# de-de sheet: new row 4
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/60aed80794e1226e8469026e0ce8cf0cadddc531/e2e/$mdName",
    $null,
    $null,
    $mdName)
$wsDeDe.Range("B4").Value = $statusInSync
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/acc7b5926c34f00e05dc30f00421dc3f3fc2e5cf/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlfName",
    $null,
    $null,
    $deXlfName)
$wsDeDe.Range("D4").Value = "2016-01-25 10:47:48"
$wsDeDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/ab280a8d3807a70e01ea14ab59647c8c58f4016f/e2e/$mdName",
    $null,
    $null,
    $mdName)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ab5eba3438fd7c47d89ac5f3b5d0f4212e98368b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlfName",
    $null,
    $null,
    $deXlfName)
$wsDeDe.Range("G4").Value = "2016-01-25 10:48:36"
$wsDeDe.Range("H4").Value = "Include"
